$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Abhishek Sharma"

# Insert a new first column (matchNo) - shifts existing teamName..result from A:L to B:M
$ws.Columns.Item(1).Insert()

# Insert three new data rows above the existing data row (old row 2 -> becomes row 5)
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

# Insert three new data rows below the existing data row (after row 5)
$ws.Rows.Item(6).Insert()
$ws.Rows.Item(6).Insert()
$ws.Rows.Item(6).Insert()

# Header row
$headers = @("matchNo","teamName","batterName","states","runs","balls","fours","sixes","sr","opponentTeamName","venue","date","result")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = "'" + $headers[$c]
}

# All data rows (rows 2..8), values forced to text with a leading apostrophe
$data = @(
    @("44th","Sunrisers Hyderabad","Abhishek Sharma","c du Plessis b Hazlewood","18","13","1","1","138.46","Chennai Super Kings","Sharjah","September 30","Super Kings won by 6 wickets (with 2 balls remaining)"),
    @("49th","Sunrisers Hyderabad","Abhishek Sharma","st †Karthik b Shakib Al Hasan","6","10","0","0","60.00","Kolkata Knight Riders","Dubai (DSC)","October 03","KKR won by 6 wickets (with 2 balls remaining)"),
    @("52nd","Sunrisers Hyderabad","Abhishek Sharma","c Maxwell b Garton","13","10","1","1","130.00","Royal Challengers Bangalore","Abu Dhabi","October 06","Sunrisers won by 4 runs"),
    @("55th","Sunrisers Hyderabad","Abhishek Sharma","c Coulter-Nile b Neesham","33","16","4","1","206.25","Mumbai Indians","Abu Dhabi","October 08","Mumbai won by 42 runs"),
    @("20th","Sunrisers Hyderabad","Abhishek Sharma","lbw b Patel","5","6","0","0","83.33","Delhi Capitals","Chennai","April 25","Match tied (Capitals won the one-over eliminator)"),
    @("40th","Sunrisers Hyderabad","Abhishek Sharma","","21","16","1","1","131.25","Rajasthan Royals","Dubai (DSC)","September 27","Sunrisers won by 7 wickets (with 9 balls remaining)"),
    @("9th","Sunrisers Hyderabad","Abhishek Sharma","c Milne b Chahar","2","4","0","0","50.00","Mumbai Indians","Chennai","April 17","Mumbai won by 13 runs")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = "'" + $row[$c]
    }
}
